$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.443.88"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.850.70"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("D4").Value = "'0.9989"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'241.05"
$ws.Range("D6").Value = "'0.6327"
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.949.01"
$ws.Range("E8").Value = "  +107.05%  "
$ws.Range("D9").Value = "4.215.77"
$ws.Range("E9").Value = "  +94.75%  "
$ws.Range("D10").Value = "'0.07568"
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("D11").Value = "'0.2968"
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("D12").Value = "'24.66"
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("D13").Value = "'0.07728"
$ws.Range("E13").Value = "  +1.30%  "
$ws.Range("D14").Value = "'4.995"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").Value = "'0.6855"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("D16").Value = "'83.02"
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("D17").Value = "'0.000009916"
$ws.Range("E17").Value = "  +3.77%  "
$ws.Range("D18").Value = "'6.204"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").Value = "29.478.86"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").Value = "'231.91"
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("D22").Value = "'0.9998"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  -1.87%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "'155.76"
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("E26").Value = "  -1.43%  "
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("E28").Value = "  -0.54%  "
$ws.Range("D29").Value = "4.213.74"
$ws.Range("E29").Value = "  +104.91%  "
$ws.Range("D30").Value = "'1.469"
$ws.Range("E30").Value = "  -1.39%  "
$ws.Range("D31").Value = "'0.05804"
$ws.Range("E31").Value = "  -3.18%  "
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("D33").Value = "'4.134"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("D34").Value = "'4.023"
$ws.Range("E34").Value = "  -1.23%  "
$ws.Range("D35").Value = "'1.860"
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("D36").Value = "'1.160"
$ws.Range("E36").Value = "  -1.33%  "
$ws.Range("D37").Value = "'0.7172"
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("D38").Value = "'2.595"
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("D39").Value = "1.251.64"
$ws.Range("E39").Value = "  +3.97%  "
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("E41").Value = "  +1.61%  "
$ws.Range("D42").Value = "'0.9019"
$ws.Range("E42").Value = "  -0.96%  "
$ws.Range("D43").Value = "'6.116"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "'101.75"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").Value = "'67.11"
$ws.Range("E46").Value = "  +0.43%  "
$ws.Range("D47").Value = "'7.207"
$ws.Range("E47").Value = "  -1.15%  "
$ws.Range("D48").Value = "'9.148"
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("D49").Value = "'0.4023"
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("D50").Value = "'1.687"
$ws.Range("E50").Value = "  +1.68%  "
$ws.Range("D51").Value = "'0.1127"
$ws.Range("E51").Value = "  -0.11%  "
